$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (08-nov) before the
#     existing "01-oct." column (column DM), shifting DM:EQ -> DN:ER.
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("DM").Insert()

$ws1.Range("DM1").Value = "08-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 117).Value = "-"
}

# --- Sheet "Gaz": append next day's row.
#     Leading apostrophe forces text (avoids auto date-serial conversion);
#     Style reset to "Normal" drops the date number-format Excel would
#     otherwise stamp on the cell.
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A145").Value = "'2025-11-06"
$ws2.Range("A145").Style = "Normal"
$ws2.Range("B145").Value = 30.35

# --- Sheet "CO2": append next day's row.
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A145").Value = "'2025-11-06"
$ws3.Range("A145").Style = "Normal"
$ws3.Range("B145").Value = 79.94
